# Generate Report for Handoff
# Regenerates the handoff report: the source file is identified by a new
# GUID-named markdown file and a new handoff commit/package hash, and the
# handoff timestamps for each locale are bumped to the new run's time.

$wb = $excel.ActiveWorkbook

$newMd = "3a04710b-64ed-4dc4-8951-53956a60b171.md"
$newZh = "3a04710b-64ed-4dc4-8951-53956a60b171.2b1597d1417fb273211f03b22390d9483163e21f.zh-cn.xlf"
$newDe = "3a04710b-64ed-4dc4-8951-53956a60b171.2b1597d1417fb273211f03b22390d9483163e21f.de-de.xlf"

$newZhTime = "2016-02-25 06:08:33"
$newDeTime = "2016-02-25 06:08:44"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- update cell values (Source File Name / Latest Handoff File / Latest Handoff Datetime) ---

$wsOverview.Range("A2").Value2 = $newMd

$wsZh.Range("A2").Value2 = $newMd
$wsZh.Range("C2").Value2 = $newZh
$wsZh.Range("D2").Value2 = $newZhTime

$wsDe.Range("A2").Value2 = $newMd
$wsDe.Range("C2").Value2 = $newDe
$wsDe.Range("D2").Value2 = $newDeTime

# --- update hyperlink display text to match the new file names ---

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    }
}

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.Range.Address() -eq '$C$2') {
        $hl.TextToDisplay = $newZh
    }
}

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.Range.Address() -eq '$C$2') {
        $hl.TextToDisplay = $newDe
    }
}
